# Word COM-interop script (PowerShell-style) implementing the commit:
# "Updated abstracts for several courses."
#
# Two content-level edits are applied to word/document.xml:
#   1. Remove the stray `_GoBack` bookmark that Word leaves behind after
#      the last edit position (bookmarkStart/bookmarkEnd pair right after
#      the "ERC 1.1" run).
#   2. Trim the "For more information" section: drop the explicit list of
#      social-media / news links (IBM Training News, YouTube, Facebook,
#      Twitter) together with the intro sentence's text, collapsing the
#      five paragraphs into a single empty paragraph that keeps the
#      left-indent formatting (0.5in / 720 twips) that used to belong to
#      the link list.

$d = $word.ActiveDocument

# --- 1. Remove the leftover "_GoBack" bookmark -----------------------------
$goBack = $d.Bookmarks("_GoBack")
$goBack.Delete()

# --- 2. Collapse the "stay informed" link list -----------------------------
# Locate the paragraph that introduces the link list and the following
# paragraphs that each hold one social/news link, by matching on text.
$paras = $d.Paragraphs
$introIndex = 0
$lastLinkIndex = 0
$i = 0
foreach ($p in $paras) {
    $i = $i + 1
    $t = $p.Range.Text
    if ($t.StartsWith("To stay informed about IBM training")) {
        $introIndex = $i
    }
    if ($t.StartsWith("Twitter: twitter.com/IBMCloudEdu")) {
        $lastLinkIndex = $i
    }
}

if ($introIndex -gt 0 -and $lastLinkIndex -gt $introIndex) {
    # Delete the four link paragraphs entirely (IBM Training News / YouTube /
    # Facebook / Twitter), paragraph marks included, leaving the intro
    # paragraph's own mark untouched.
    $firstLinkPara = $paras.Item($introIndex + 1)
    $lastLinkPara = $paras.Item($lastLinkIndex)
    $deleteRange = $d.Range($firstLinkPara.Range.Start, $lastLinkPara.Range.End)
    $deleteRange.Delete()

    # Clear the intro sentence's own text, but keep its paragraph mark so a
    # single empty paragraph remains.
    $introPara = $paras.Item($introIndex)
    $introTextRange = $d.Range($introPara.Range.Start, $introPara.Range.End - 1)
    $introTextRange.Text = ""

    # The now-empty paragraph inherits the 0.5in left indent that used to
    # apply to the link paragraphs.
    $paras.Item($introIndex).Format.LeftIndent = 36
}
